$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 4185.7144
$ws.Range("I29").Value = 2150
$ws.Range("K29").Value = 6450
$ws.Range("M29").Value = -6169
$ws.Range("H62").Value = 3210.3635
$ws.Range("I62").Value = 3210.3635
$ws.Range("K62").Value = 3210.3635
$ws.Range("M62").Value = -2586.3635
$ws.Range("H65").Value = 3210.3635
$ws.Range("I65").Value = 3210.3635
$ws.Range("K65").Value = 16051.8175
$ws.Range("M65").Value = -12931.8175
$ws.Range("H98").Value = 25919.87
$ws.Range("I98").Value = 28521.088
$ws.Range("K98").Value = 28521.088
$ws.Range("M98").Value = -27023.088
$ws.Range("H116").Value = 6101.7856
$ws.Range("I116").Value = 5779.5
$ws.Range("K116").Value = 5779.5
$ws.Range("M116").Value = -2337.5
$ws.Range("H122").Value = 25919.87
$ws.Range("I122").Value = 28521.088
$ws.Range("K122").Value = 85563.264
$ws.Range("M122").Value = -83113.264
$ws.Range("H127").Value = 4200.3076
$ws.Range("I127").Value = 4230.625
$ws.Range("J127").Value = 4151.8
$ws.Range("K127").Value = 12691.875
$ws.Range("L127").Value = 12455.4
$ws.Range("M127").Value = -7731.875
$ws.Range("N127").Value = -22375.4
$ws.Range("H132").Value = 2023.2565
$ws.Range("I132").Value = 2126.7576
$ws.Range("K132").Value = 6380.2728
$ws.Range("M132").Value = -3850.2728
$ws.Range("H138").Value = 4064
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4064
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 12192
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -22472

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 7979
$ws.Range("I25").Value = 3947.5
$ws.Range("J25").Value = 10666.667
$ws.Range("K25").Value = 3947.5
$ws.Range("L25").Value = 10666.667
$ws.Range("M25").Value = -3545.5
$ws.Range("N25").Value = -11470.667
$ws.Range("H32").Value = 3872.44
$ws.Range("I32").Value = 2976.806
$ws.Range("K32").Value = 2976.806
$ws.Range("M32").Value = -2689.806
$ws.Range("H35").Value = 6787.25
$ws.Range("J35").Value = 7038.8
$ws.Range("L35").Value = 7038.8
$ws.Range("N35").Value = -7850.8
$ws.Range("H74").Value = 493225.22
$ws.Range("I74").Value = 840062.1
$ws.Range("K74").Value = 840062.1
$ws.Range("M74").Value = -839188.1
$ws.Range("H77").Value = 493225.22
$ws.Range("I77").Value = 840062.1
$ws.Range("K77").Value = 4200310.5
$ws.Range("M77").Value = -4195942.5
$ws.Range("H132").Value = 1877.8524
$ws.Range("I132").Value = 1810.585
$ws.Range("K132").Value = 5431.755
$ws.Range("M132").Value = -2901.755

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 955.3333
$ws.Range("J20").Value = 803.5454999999999
$ws.Range("L20").Value = 803.5454999999999
$ws.Range("N20").Value = -1297.5455

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1380.2034
$ws.Range("I31").Value = 1296.3469
$ws.Range("K31").Value = 1296.3469
$ws.Range("M31").Value = -1001.3469
$ws.Range("H34").Value = 1380.2034
$ws.Range("I34").Value = 1296.3469
$ws.Range("K34").Value = 1296.3469
$ws.Range("M34").Value = -1094.3469
$ws.Range("H58").Value = 2733.1904
$ws.Range("I58").Value = 2322.5
$ws.Range("J58").Value = 3106.5454
$ws.Range("K58").Value = 2322.5
$ws.Range("L58").Value = 3106.5454
$ws.Range("M58").Value = -2119.5
$ws.Range("N58").Value = -3512.5454
$ws.Range("H99").Value = 4618.4614
$ws.Range("I99").Value = 4618.4614
$ws.Range("K99").Value = 4618.4614
$ws.Range("M99").Value = -3120.4614
$ws.Range("H107").Value = 25665584
$ws.Range("I107").Value = 40036624
$ws.Range("J107").Value = 3009.2856
$ws.Range("K107").Value = 40036624
$ws.Range("L107").Value = 3009.2856
$ws.Range("M107").Value = -40034704
$ws.Range("N107").Value = -6849.2856
$ws.Range("H126").Value = 4618.4614
$ws.Range("I126").Value = 4618.4614
$ws.Range("K126").Value = 13855.3842
$ws.Range("M126").Value = -11385.3842
$ws.Range("H136").Value = 2733.1904
$ws.Range("I136").Value = 2322.5
$ws.Range("J136").Value = 3106.5454
$ws.Range("K136").Value = 6967.5
$ws.Range("L136").Value = 9319.636200000001
$ws.Range("M136").Value = -4417.5
$ws.Range("N136").Value = -14419.6362

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 202
$ws.Range("I86").Value = 202
$ws.Range("K86").Value = 606
$ws.Range("M86").Value = 580
$ws.Range("H89").Value = 202
$ws.Range("I89").Value = 202
$ws.Range("K89").Value = 1818
$ws.Range("M89").Value = 4110
$ws.Range("H113").Value = 2166
$ws.Range("J113").Value = 3365.1667
$ws.Range("L113").Value = 10095.5001
$ws.Range("N113").Value = -14435.5001
$ws.Range("H122").Value = 345.66666
$ws.Range("J122").Value = 345.66666
$ws.Range("L122").Value = 3110.99994
$ws.Range("N122").Value = -8010.99994
$ws.Range("H137").Value = 2035.9436
$ws.Range("I137").Value = 1499
$ws.Range("K137").Value = 4497
$ws.Range("M137").Value = 603

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2250.5
$ws.Range("I80").Value = 2124.5
$ws.Range("J80").Value = 2502.5
$ws.Range("K80").Value = 2124.5
$ws.Range("L80").Value = 2502.5
$ws.Range("M80").Value = -1126.5
$ws.Range("N80").Value = -4498.5
$ws.Range("H83").Value = 2250.5
$ws.Range("I83").Value = 2124.5
$ws.Range("J83").Value = 2502.5
$ws.Range("K83").Value = 10622.5
$ws.Range("L83").Value = 12512.5
$ws.Range("M83").Value = -5630.5
$ws.Range("N83").Value = -22496.5
$ws.Range("H107").Value = 1039.5333
$ws.Range("J107").Value = 1554.5
$ws.Range("L107").Value = 1554.5
$ws.Range("N107").Value = -5394.5
$ws.Range("H132").Value = 2419.62
$ws.Range("I132").Value = 2387.125
$ws.Range("J132").Value = 3199.5
$ws.Range("K132").Value = 7161.375
$ws.Range("L132").Value = 9598.5
$ws.Range("M132").Value = -4631.375
$ws.Range("N132").Value = -14658.5

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 39506.582
$ws.Range("I7").Value = 53312.25
$ws.Range("K7").Value = 53312.25
$ws.Range("M7").Value = -53200.25
$ws.Range("H40").Value = 26316.588
$ws.Range("I40").Value = 27867.812
$ws.Range("K40").Value = 27867.812
$ws.Range("M40").Value = -27731.812
$ws.Range("H55").Value = 1040.1
$ws.Range("J55").Value = 1979.8889
$ws.Range("L55").Value = 1979.8889
$ws.Range("N55").Value = -2325.8889
$ws.Range("H61").Value = 71436310
$ws.Range("J61").Value = 10509.25
$ws.Range("L61").Value = 10509.25
$ws.Range("N61").Value = -10913.25
$ws.Range("H113").Value = 71436310
$ws.Range("J113").Value = 10509.25
$ws.Range("L113").Value = 10509.25
$ws.Range("N113").Value = -14849.25
$ws.Range("H126").Value = 39506.582
$ws.Range("I126").Value = 53312.25
$ws.Range("K126").Value = 159936.75
$ws.Range("M126").Value = -157466.75
$ws.Range("H132").Value = 2455.0862
$ws.Range("I132").Value = 2165.3022
$ws.Range("J132").Value = 3285.8
$ws.Range("K132").Value = 6495.9066
$ws.Range("L132").Value = 9857.400000000001
$ws.Range("M132").Value = -3965.9066
$ws.Range("N132").Value = -14917.4

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 9415.666999999999
$ws.Range("I28").Value = 7874.75
$ws.Range("K28").Value = 7874.75
$ws.Range("M28").Value = -7526.75
$ws.Range("H41").Value = 20067.625
$ws.Range("I41").Value = 17999.5
$ws.Range("J41").Value = 20757
$ws.Range("K41").Value = 17999.5
$ws.Range("L41").Value = 20757
$ws.Range("M41").Value = -17609.5
$ws.Range("N41").Value = -21537
$ws.Range("H122").Value = 2464.3428
$ws.Range("I122").Value = 2317.0667
$ws.Range("K122").Value = 6951.2001
$ws.Range("M122").Value = -4501.2001
$ws.Range("H126").Value = 30307064
$ws.Range("I126").Value = 34486680
$ws.Range("K126").Value = 103460040
$ws.Range("M126").Value = -103457570
$ws.Range("H132").Value = 1989.1578
$ws.Range("I132").Value = 1858.898
$ws.Range("J132").Value = 2787
$ws.Range("K132").Value = 5576.694
$ws.Range("L132").Value = 8361
$ws.Range("M132").Value = -3046.694
$ws.Range("N132").Value = -13421
